$d = $word.ActiveDocument
$p1 = $d.Paragraphs.Item(1)

# Add a paragraph border (no line drawn, just 5-twip padding on all 4 sides)
# to the first paragraph's pPr (-> <w:pBdr><w:top w:space="5"/>... />).
$borders = $p1.Borders
$borders.DistanceFromTop = 5
$borders.DistanceFromLeft = 5
$borders.DistanceFromBottom = 5
$borders.DistanceFromRight = 5

# Widen the left indent of paragraph 1 from 120 twips (6pt) to 225 twips (11.25pt).
$p1.Range.ParagraphFormat.LeftIndent = 11.25

# The paragraph used to be "**ID__AFFARS_5332_topic_2__ID**" followed by a
# separate run containing a single trailing space. Drop that trailing-space
# run (it sits right before the paragraph mark).
$pEnd = $p1.Range.End
$spaceRange = $d.Range($pEnd - 2, $pEnd - 1)
if ($spaceRange.Text -eq " ") {
    $spaceRange.Delete()
}

# Rename the placeholder token itself.
$d.Content.Find.Execute("**ID__AFFARS_5332_topic_2__ID**", $true, $false, $false, $false, $false, $true, 1, $false, "**ID__AFFARS_SUBPART_5332_1__ID**", 2)
